$wb = $excel.ActiveWorkbook

$oldGuid = "4b4c19a2-07bd-4260-a3ed-2607b79109fd"
$newGuid = "1b7b9085-2de4-499a-bd9d-243bea44e46e"
$oldHash = "1abd6df45ed39ce7c1f45640f125c11961bf1118"
$newHash = "193503295b4c241d68c95668603e96e7c9b9a84b"

# ---------------------------------------------------------------------------
# Overview sheet - refresh the generated file name / path and the HO xliff
# generation timestamp. The existing B2 hyperlink keeps its original target,
# only the displayed text changes (matches the new file name).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = ($newGuid + ".md")
$wsOverview.Range("G2").Value = "2016-09-02 01:11:13"

$overviewLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec0c1eb1c008354f9286e488c21f32ca8fad9b/e2e/" + $oldGuid + ".md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddress, [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $newGuid + ".md"))

# ---------------------------------------------------------------------------
# zh-cn sheet - new source file name, new handoff xliff name + handoff time,
# handback is no longer available (target/handback file cleared, handback
# datetime reset to the zero date), and the old "Latest Target File" hyperlink
# on I2 is removed along with its now-empty cell's hyperlink styling.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = ($newGuid + ".md")
$wsZh.Range("G2").Value = ($newGuid + "." + $newHash + ".zh-cn.xlf")
$wsZh.Range("H2").Value = "2016-09-02 01:11:08"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Style = "Normal"

$zhLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec0c1eb1c008354f9286e488c21f32ca8fad9b/e2e/" + $oldGuid + ".md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhLinkAddress, [System.Type]::Missing, [System.Type]::Missing, ($newGuid + ".md"))

# ---------------------------------------------------------------------------
# de-de sheet - same shape of change as zh-cn.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = ($newGuid + ".md")
$wsDe.Range("G2").Value = ($newGuid + "." + $newHash + ".de-de.xlf")
$wsDe.Range("H2").Value = "2016-09-02 01:11:13"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Style = "Normal"

$deLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec0c1eb1c008354f9286e488c21f32ca8fad9b/e2e/" + $oldGuid + ".md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deLinkAddress, [System.Type]::Missing, [System.Type]::Missing, ($newGuid + ".md"))

# ---------------------------------------------------------------------------
# Column widths on zh-cn / de-de: "Latest Target File" (I) and
# "Latest Handback File" (J) shrink now that they hold no long file names.
# ---------------------------------------------------------------------------
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
